$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 267.93332
$ws.Range("I33").Value = 243.54546
$ws.Range("K33").Value = 243.54546
$ws.Range("M33").Value = -14.54545999999999
$ws.Range("H55").Value = 479.4
$ws.Range("I55").Value = 592
$ws.Range("K55").Value = 592
$ws.Range("M55").Value = -378
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("H106").Value = 24842.375
$ws.Range("I106").Value = 24842.375
$ws.Range("K106").Value = 24842.375
$ws.Range("M106").Value = -24211.375
$ws.Range("H107").Value = 2373
$ws.Range("I107").Value = 2444.375
$ws.Range("K107").Value = 2444.375
$ws.Range("M107").Value = -524.375
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H137").Value = 2005.5
$ws.Range("J137").Value = 2019.8
$ws.Range("L137").Value = 6059.4
$ws.Range("N137").Value = -11159.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6133.0347
$ws.Range("I32").Value = 4660.2915
$ws.Range("J32").Value = 13202.2
$ws.Range("K32").Value = 4660.2915
$ws.Range("L32").Value = 13202.2
$ws.Range("M32").Value = -4373.2915
$ws.Range("N32").Value = -13776.2
$ws.Range("H61").Value = 7949.25
$ws.Range("I61").Value = 7949.25
$ws.Range("K61").Value = 7949.25
$ws.Range("M61").Value = -7737.25
$ws.Range("H74").Value = 1114.1
$ws.Range("I74").Value = 955.25
$ws.Range("J74").Value = 1749.5
$ws.Range("K74").Value = 955.25
$ws.Range("L74").Value = 1749.5
$ws.Range("M74").Value = -81.25
$ws.Range("N74").Value = -3497.5
$ws.Range("H77").Value = 1114.1
$ws.Range("I77").Value = 955.25
$ws.Range("J77").Value = 1749.5
$ws.Range("K77").Value = 4776.25
$ws.Range("L77").Value = 8747.5
$ws.Range("M77").Value = -408.25
$ws.Range("N77").Value = -17483.5
$ws.Range("H122").Value = 784393.3
$ws.Range("I122").Value = 1119354.8
$ws.Range("J122").Value = 30730
$ws.Range("K122").Value = 3358064.4
$ws.Range("L122").Value = 92190
$ws.Range("M122").Value = -3355614.4
$ws.Range("N122").Value = -97090
$ws.Range("H132").Value = 4142.7144
$ws.Range("I132").Value = 2999.5
$ws.Range("K132").Value = 8998.5
$ws.Range("M132").Value = -6468.5
$ws.Range("H136").Value = 7949.25
$ws.Range("I136").Value = 7949.25
$ws.Range("K136").Value = 23847.75
$ws.Range("M136").Value = -21297.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 49166.668
$ws.Range("I48").Value = 46000
$ws.Range("J48").Value = 65000
$ws.Range("K48").Value = 46000
$ws.Range("L48").Value = 65000
$ws.Range("M48").Value = -45524
$ws.Range("N48").Value = -65952
$ws.Range("H58").Value = 2931.8462
$ws.Range("I58").Value = 1681.2222
$ws.Range("J58").Value = 5745.75
$ws.Range("K58").Value = 1681.2222
$ws.Range("L58").Value = 5745.75
$ws.Range("M58").Value = -1478.2222
$ws.Range("N58").Value = -6151.75
$ws.Range("H136").Value = 2931.8462
$ws.Range("I136").Value = 1681.2222
$ws.Range("J136").Value = 5745.75
$ws.Range("K136").Value = 5043.6666
$ws.Range("L136").Value = 17237.25
$ws.Range("M136").Value = -2493.6666
$ws.Range("N136").Value = -22337.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 419.75
$ws.Range("J2").Value = 545.3333
$ws.Range("L2").Value = 3271.9998
$ws.Range("N2").Value = -3497.9998
$ws.Range("H33").Value = 142929.86
$ws.Range("I33").Value = 77.5
$ws.Range("K33").Value = 465
$ws.Range("M33").Value = -182
$ws.Range("H131").Value = 1714.8036
$ws.Range("I131").Value = 935
$ws.Range("J131").Value = 1774.7885
$ws.Range("K131").Value = 2805
$ws.Range("L131").Value = 5324.3655
$ws.Range("M131").Value = 2235
$ws.Range("N131").Value = -15404.3655
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2716
$ws.Range("I97").Value = 2199.875
$ws.Range("J97").Value = 3091.3635
$ws.Range("K97").Value = 2199.875
$ws.Range("L97").Value = 3091.3635
$ws.Range("M97").Value = -1703.875
$ws.Range("N97").Value = -4083.3635
$ws.Range("H122").Value = 203199.2
$ws.Range("H132").Value = 1683.6154
$ws.Range("I132").Value = 1444.3636
$ws.Range("K132").Value = 4333.0908
$ws.Range("M132").Value = -1803.0908
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1099.4615
$ws.Range("I16").Value = 1099.4615
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1099.4615
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -929.4614999999999
$ws.Range("N16").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 63500
$ws.Range("J64").Value = 63500
$ws.Range("L64").Value = 63500
$ws.Range("N64").Value = -63996
$ws.Range("H67").Value = 63500
$ws.Range("J67").Value = 63500
$ws.Range("L67").Value = 63500
$ws.Range("N67").Value = -65216
$ws.Range("H74").Value = 60000
$ws.Range("J74").Value = 60000
$ws.Range("L74").Value = 60000
$ws.Range("N74").Value = -61872
$ws.Range("H77").Value = 60000
$ws.Range("J77").Value = 60000
$ws.Range("L77").Value = 180000
$ws.Range("N77").Value = -189360
$ws.Range("H107").Value = 874.875
$ws.Range("I107").Value = 714.1429000000001
$ws.Range("K107").Value = 2142.4287
$ws.Range("M107").Value = -222.4287000000004
$ws.Range("H122").Value = 2581.625
$ws.Range("I122").Value = 2615.4285
$ws.Range("J122").Value = 2345
$ws.Range("K122").Value = 7846.2855
$ws.Range("L122").Value = 7035
$ws.Range("M122").Value = -5396.2855
$ws.Range("N122").Value = -11935
$ws.Range("H136").Value = 852.5
$ws.Range("I136").Value = 686.8421
$ws.Range("K136").Value = 2060.5263
$ws.Range("M136").Value = 489.4737
